# Generate Report for Handoff
# Refresh the "Latest Handoff Date(time)" for every file that is not yet
# fully handed back (i.e. still "Handback transform failed",
# "In Translation" stays untouched, or "Ready for handoff") to the new
# report-generation timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value  = "2016-28-13 16:28:17"
$overview.Range("D9").Value  = "2016-28-13 16:28:17"
$overview.Range("D10").Value = "2016-28-13 16:28:17"
$overview.Range("D11").Value = "2016-28-13 16:28:17"
$overview.Range("D12").Value = "2016-28-13 16:28:17"
$overview.Range("D13").Value = "2016-28-13 16:28:17"
$overview.Range("D14").Value = "2016-28-13 16:28:17"
$overview.Range("D15").Value = "2016-28-13 16:28:17"
$overview.Range("D16").Value = "2016-28-13 16:28:17"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value  = "2016-03-13 16:28:10"
$zhcn.Range("E9").Value  = "2016-03-13 16:28:10"
$zhcn.Range("E10").Value = "2016-03-13 16:28:10"
$zhcn.Range("E11").Value = "2016-03-13 16:28:10"
$zhcn.Range("E12").Value = "2016-03-13 16:28:10"
$zhcn.Range("E13").Value = "2016-03-13 16:28:10"
$zhcn.Range("E14").Value = "2016-03-13 16:28:10"
$zhcn.Range("E15").Value = "2016-03-13 16:28:10"
$zhcn.Range("E16").Value = "2016-03-13 16:28:10"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value  = "2016-03-13 16:28:17"
$dede.Range("E9").Value  = "2016-03-13 16:28:17"
$dede.Range("E10").Value = "2016-03-13 16:28:17"
$dede.Range("E11").Value = "2016-03-13 16:28:17"
$dede.Range("E12").Value = "2016-03-13 16:28:17"
$dede.Range("E13").Value = "2016-03-13 16:28:17"
$dede.Range("E14").Value = "2016-03-13 16:28:17"
$dede.Range("E15").Value = "2016-03-13 16:28:17"
$dede.Range("E16").Value = "2016-03-13 16:28:17"
